$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 5) appended below the existing data (rows 1-4).
# Columns A-D are text (Date/Time/Weekday/Week) like the existing rows;
# force text storage for the date- and number-looking values so Excel
# doesn't auto-convert them, then clear the temporary "@" number format
# so no extra style is left behind on the cells.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"

$ws.Range("A5").Value = "2023-05-28"
$ws.Range("B5").Value = "22:07:21"
$ws.Range("C5").Value = "Sunday"
$ws.Range("D5").Value = "21"

$ws.Range("A5").ClearFormats()
$ws.Range("D5").ClearFormats()

# Columns E-T are the numeric per-city resale figures.
$ws.Range("E5").Value = 119957
$ws.Range("F5").Value = 133322
$ws.Range("G5").Value = 157974
$ws.Range("H5").Value = 130719
$ws.Range("I5").Value = 174322
$ws.Range("J5").Value = 114355
$ws.Range("K5").Value = 198289
$ws.Range("L5").Value = 219973
$ws.Range("M5").Value = 172034
$ws.Range("N5").Value = 119759
$ws.Range("O5").Value = 38627
$ws.Range("P5").Value = 34955
$ws.Range("Q5").Value = 50395
$ws.Range("R5").Value = -1
$ws.Range("S5").Value = 37182
$ws.Range("T5").Value = -1
